$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Find the last used row in column A (currently row 248, date 01-10-2021)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$dates = @("02-10-2021", "03-10-2021", "04-10-2021", "05-10-2021", "06-10-2021")

$r = $lastRow
foreach ($d in $dates) {
    $r = $r + 1
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $d
    $cellA.Style = "Normal"
    $ws.Cells.Item($r, 2).Value = 322
    $ws.Cells.Item($r, 3).Value = 0
}
